$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume update (GitHub Actions scheduled refresh)
# Ensure text-formatted columns (Price / Volume) keep their values as text,
# matching the original inlineStr cell type, instead of being auto-converted
# to numbers by Excel when the string looks like a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.273.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.176.29"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.63%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.98"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.610"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "69.56"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.63%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.576"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.54"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0920"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.17"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.76%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.70"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.502.57"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.98"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.181.56"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.795"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.093.02"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.91%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.61"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.90"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.45%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "225.91"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.37%  "

$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.47"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -6.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.90"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -7.57%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.71"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -7.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.52"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.59%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.14"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.92"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.50"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0763"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.08"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -10.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.120"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.101"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -7.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.04"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0283"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.07"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.58"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -12.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.41"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "59.12"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -9.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.189"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0972"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.25"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.51"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.08"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.22%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.19"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -7.41%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.58%  "
